$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 359.1111
$ws.Range("I19").Value = 375
$ws.Range("J19").Value = 346.4
$ws.Range("K19").Value = 375
$ws.Range("L19").Value = 346.4
$ws.Range("M19").Value = -200
$ws.Range("N19").Value = -696.4
$ws.Range("H68").Value = 36629.332
$ws.Range("J68").Value = 36629.332
$ws.Range("L68").Value = 36629.332
$ws.Range("N68").Value = -38127.332
$ws.Range("H71").Value = 36629.332
$ws.Range("J71").Value = 36629.332
$ws.Range("L71").Value = 109887.996
$ws.Range("N71").Value = -117375.996
$ws.Range("H93").Value = 81300.75
$ws.Range("J93").Value = 81300.75
$ws.Range("L93").Value = 81300.75
$ws.Range("N93").Value = -86292.75
$ws.Range("H111").Value = 22225376
$ws.Range("I111").Value = 30304786
$ws.Range("J111").Value = 7000
$ws.Range("K111").Value = 90914358
$ws.Range("L111").Value = 21000
$ws.Range("M111").Value = -90911291
$ws.Range("N111").Value = -27134
$ws.Range("H112").Value = 2156.976
$ws.Range("I112").Value = 300.55554
$ws.Range("J112").Value = 2663.2727
$ws.Range("K112").Value = 901.66662
$ws.Range("L112").Value = 7989.8181
$ws.Range("M112").Value = 206.33338
$ws.Range("N112").Value = -10205.8181
$ws.Range("H137").Value = 3659758
$ws.Range("I137").Value = 1786952.9
$ws.Range("J137").Value = 7693492.5
$ws.Range("K137").Value = 5360858.699999999
$ws.Range("L137").Value = 23080477.5
$ws.Range("M137").Value = -5358308.699999999
$ws.Range("N137").Value = -23085577.5
$ws.Range("H138").Value = 2291.3672
$ws.Range("I138").Value = 2480.0833
$ws.Range("J138").Value = 2257.5671
$ws.Range("K138").Value = 7440.249899999999
$ws.Range("L138").Value = 6772.701300000001
$ws.Range("M138").Value = -2300.249899999999
$ws.Range("N138").Value = -17052.7013
$ws.Range("H141").Value = 1636.8125
$ws.Range("I141").Value = 1162.5555
$ws.Range("J141").Value = 2246.5715
$ws.Range("K141").Value = 3487.6665
$ws.Range("L141").Value = 6739.7145
$ws.Range("M141").Value = 1692.3335
$ws.Range("N141").Value = -17099.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 166668140
$ws.Range("I102").Value = 1405
$ws.Range("J102").Value = 250001500
$ws.Range("K102").Value = 1405
$ws.Range("L102").Value = 250001500
$ws.Range("M102").Value = 217
$ws.Range("N102").Value = -250004744
$ws.Range("H132").Value = 98282
$ws.Range("I132").Value = 120429.79
$ws.Range("J132").Value = 5261.3
$ws.Range("K132").Value = 361289.37
$ws.Range("L132").Value = 15783.9
$ws.Range("M132").Value = -358759.37
$ws.Range("N132").Value = -20843.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4409
$ws.Range("I105").Value = 2357
$ws.Range("J105").Value = 8000
$ws.Range("K105").Value = 2357
$ws.Range("L105").Value = 8000
$ws.Range("M105").Value = -610
$ws.Range("N105").Value = -11494
$ws.Range("H132").Value = 39800
$ws.Range("J132").Value = 39800
$ws.Range("L132").Value = 39800
$ws.Range("N132").Value = -49920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 253.15384
$ws.Range("I22").Value = 244.63637
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 244.63637
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 105.36363
$ws.Range("N22").Value = -1000
$ws.Range("H58").Value = 816.11536
$ws.Range("I58").Value = 844.45
$ws.Range("J58").Value = 721.6667
$ws.Range("K58").Value = 844.45
$ws.Range("L58").Value = 721.6667
$ws.Range("M58").Value = -641.45
$ws.Range("N58").Value = -1127.6667
$ws.Range("H132").Value = 5080.0713
$ws.Range("I132").Value = 4343.3335
$ws.Range("J132").Value = 9500.5
$ws.Range("K132").Value = 13030.0005
$ws.Range("L132").Value = 28501.5
$ws.Range("M132").Value = -10500.0005
$ws.Range("N132").Value = -33561.5
$ws.Range("H136").Value = 816.11536
$ws.Range("I136").Value = 844.45
$ws.Range("J136").Value = 721.6667
$ws.Range("K136").Value = 2533.35
$ws.Range("L136").Value = 2165.0001
$ws.Range("M136").Value = 16.64999999999964
$ws.Range("N136").Value = -7265.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2548.75
$ws.Range("I122").Value = 2380
$ws.Range("J122").Value = 2830
$ws.Range("K122").Value = 7140
$ws.Range("L122").Value = 8490
$ws.Range("M122").Value = -4690
$ws.Range("N122").Value = -13390

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4764350
$ws.Range("I7").Value = 12501726
$ws.Range("J7").Value = 2888.2307
$ws.Range("K7").Value = 12501726
$ws.Range("L7").Value = 2888.2307
$ws.Range("M7").Value = -12501614
$ws.Range("N7").Value = -3112.2307
$ws.Range("H40").Value = 2154.182
$ws.Range("I40").Value = 1783.3334
$ws.Range("J40").Value = 2599.2
$ws.Range("K40").Value = 1783.3334
$ws.Range("L40").Value = 2599.2
$ws.Range("M40").Value = -1647.3334
$ws.Range("N40").Value = -2871.2
$ws.Range("H126").Value = 4764350
$ws.Range("I126").Value = 12501726
$ws.Range("J126").Value = 2888.2307
$ws.Range("K126").Value = 37505178
$ws.Range("L126").Value = 8664.6921
$ws.Range("M126").Value = -37502708
$ws.Range("N126").Value = -13604.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 27546
$ws.Range("J121").Value = 27546
$ws.Range("L121").Value = 27546
$ws.Range("N121").Value = -31040
$ws.Range("H126").Value = 2452.3447
$ws.Range("I126").Value = 1996.7368
$ws.Range("J126").Value = 3318
$ws.Range("K126").Value = 5990.2104
$ws.Range("L126").Value = 9954
$ws.Range("M126").Value = -3520.2104
$ws.Range("N126").Value = -14894
$ws.Range("H132").Value = 3170.0605
$ws.Range("I132").Value = 3456.3076
$ws.Range("J132").Value = 2106.8572
$ws.Range("K132").Value = 10368.9228
$ws.Range("L132").Value = 6320.571599999999
$ws.Range("M132").Value = -7838.9228
$ws.Range("N132").Value = -11380.5716
$ws.Range("H136").Value = 1538.1333
$ws.Range("I136").Value = 1493.8302
$ws.Range("K136").Value = 4481.4906
$ws.Range("M136").Value = -1931.4906
$ws.Range("H140").Value = 27374.889
$ws.Range("J140").Value = 27374.889
$ws.Range("L140").Value = 27374.889
$ws.Range("N140").Value = -37734.889
$ws.Range("H141").Value = 63691.25
$ws.Range("J141").Value = 63691.25
$ws.Range("L141").Value = 63691.25
$ws.Range("N141").Value = -74051.25
